# Apply the "data refresh" edit described by the diff:
#  - 展览 (Exhibition) sheet: update "想去人数" (F) counts for several rows
#  - 演出 (Performance) sheet: remove the finished "苗阜王声" event (row 2),
#       shift everything up, re-number the index column (A), and update a
#       handful of "想去人数" (F) counts
#  - 本地生活 (Local Life) sheet: update the "想去人数" (F) count
#  - 全部类型 (All types) sheet: update "想去人数" (F) counts for several rows

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) 展览 sheet - F column updates only
# ---------------------------------------------------------------------------
$wsExhibit = $wb.Worksheets.Item("展览")
$exhibitUpdates = @{
    2  = 2685
    3  = 582
    5  = 306
    7  = 504
    8  = 1244
    9  = 583
    10 = 316
    11 = 10
    12 = 136
    13 = 375
    14 = 5846
    15 = 96
    16 = 1812
    17 = 4280
    18 = 443
    20 = 305
    21 = 4974
    22 = 6416
    23 = 141
    25 = 705
    26 = 3829
    27 = 510
    29 = 202
    31 = 1003
    32 = 1431
    33 = 499
    34 = 598
    35 = 1621
    36 = 211
    37 = 1756
    38 = 211
    39 = 1161
    41 = 642
    42 = 99
    43 = 3486
    45 = 306
    47 = 10
    48 = 34
    49 = 3904
}
foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Range("F$row").Value = $exhibitUpdates[$row]
}

# ---------------------------------------------------------------------------
# 2) 演出 sheet - delete the finished event row, then fix up the index
#    column and a few "想去人数" counts that differ from the row that
#    shifted into their place.
# ---------------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")

# Row 2 ("杭州·苗阜王声 青曲社相声全国巡演") is gone in the refreshed export.
$wsShow.Rows.Item(2).Delete()

# After the shift the sheet now has 30 data rows (2..31). Column A is a
# simple 0-based running index (A2=1, A3=2, ... A31=30) independent of the
# row's other content, so rewrite it to match the new row positions.
for ($row = 2; $row -le 31; $row++) {
    $wsShow.Range("A$row").Value = $row - 1
}

# A few "想去人数" counts moved along with their row but the refreshed
# export shows slightly different (larger) numbers for these particular
# events - patch them explicitly.
$showUpdates = @{
    4  = 1217
    18 = 3
    22 = 1
    28 = 76
}
foreach ($row in $showUpdates.Keys) {
    $wsShow.Range("F$row").Value = $showUpdates[$row]
}

# ---------------------------------------------------------------------------
# 3) 本地生活 sheet - F column update only
# ---------------------------------------------------------------------------
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F2").Value = 4031

# ---------------------------------------------------------------------------
# 4) 全部类型 sheet - F column updates only
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$allUpdates = @{
    2  = 4031
    4  = 582
    6  = 306
    7  = 1217
    11 = 504
    13 = 1244
    14 = 583
    15 = 316
    16 = 136
    17 = 375
    18 = 96
    19 = 1812
    20 = 4280
    21 = 4974
    22 = 141
    24 = 705
    25 = 3829
    26 = 510
    28 = 202
    30 = 1431
    31 = 499
    32 = 598
    33 = 1621
    34 = 211
    35 = 1756
    38 = 642
    40 = 99
    41 = 76
    42 = 3486
    45 = 306
    47 = 10
    49 = 3904
}
foreach ($row in $allUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allUpdates[$row]
}
